$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 99, shifting existing rows 99:105 down to 100:106
$ws.Rows("99:99").Insert()

# Populate the newly inserted row 99 with the new weekly price record
$ws.Range("A99").Value = 1
$ws.Range("B99").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C99").Value = "Arica y Parinacota"
$ws.Range("D99").Value = 45106
$ws.Range("E99").Value = 15
$ws.Range("F99").Value = 100112012
$ws.Range("G99").Value = "Espinaca"
$ws.Range("H99").Value = "Sin especificar"
$ws.Range("I99").Value = "Primera"
$ws.Range("J99").Value = 250
$ws.Range("K99").Value = 800
$ws.Range("L99").Value = 1000
$ws.Range("M99").Value = 920
$ws.Range("N99").Value = "$/atado 2,5 a 3 kilos"
$ws.Range("O99").Value = "Región de Arica y Parinacota"
$ws.Range("P99").Value = 307
$ws.Range("Q99").Value = 3
$ws.Range("R99").Value = "Hortaliza"
